$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "description" column for the first block of test data
$ws.Range("C3").Value = "This is the first test"
$ws.Range("C4").Value = "This is the second test"
$ws.Range("C5").Value = "This is the third test"

# Add a new "Delete customer" test section
$ws.Range("A11").Value = "deleteCustomer"
$ws.Range("A12").Value = "id"
$ws.Range("C20").Value = " "

# Move the active selection to match the post-edit state
$ws.Range("F14").Select()
